# Auto-generated edit script applying the Moogle_Profits market-data refresh diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 6412.919
$ws.Range("J43").Value = 2895.8
$ws.Range("L43").Value = 2895.8
$ws.Range("N43").Value = -3033.8

$ws.Range("H99").Value = 1233.2307
$ws.Range("I99").Value = 419.41666
$ws.Range("J99").Value = 10999
$ws.Range("K99").Value = 1258.24998
$ws.Range("L99").Value = 32997
$ws.Range("M99").Value = 239.7500199999999
$ws.Range("N99").Value = -35993

$ws.Range("H107").Value = 2750
$ws.Range("I107").Value = 2750
$ws.Range("K107").Value = 2750
$ws.Range("M107").Value = -830

$ws.Range("H112").Value = 5449.737
$ws.Range("I112").Value = 2421.6667
$ws.Range("K112").Value = 7265.000100000001
$ws.Range("M112").Value = -6157.000100000001

$ws.Range("H135").Value = 625.9375
$ws.Range("I135").Value = 625.9375
$ws.Range("K135").Value = 5633.4375
$ws.Range("M135").Value = -3098.4375

$ws.Range("H138").Value = 2114.6938
$ws.Range("I138").Value = 2241
$ws.Range("J138").Value = 2058.9707
$ws.Range("K138").Value = 6723
$ws.Range("L138").Value = 6176.9121
$ws.Range("M138").Value = -1583
$ws.Range("N138").Value = -16456.9121

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 486.6
$ws.Range("I4").Value = 347
$ws.Range("J4").Value = 696
$ws.Range("K4").Value = 347
$ws.Range("L4").Value = 696
$ws.Range("M4").Value = -231
$ws.Range("N4").Value = -928

$ws.Range("H32").Value = 6248.409
$ws.Range("I32").Value = 2974.524
$ws.Range("J32").Value = 75000
$ws.Range("K32").Value = 2974.524
$ws.Range("L32").Value = 75000
$ws.Range("M32").Value = -2687.524
$ws.Range("N32").Value = -75574

$ws.Range("H61").Value = 4588.5
$ws.Range("I61").Value = 2537.2727
$ws.Range("K61").Value = 2537.2727
$ws.Range("M61").Value = -2325.2727

$ws.Range("H74").Value = 2390.2563
$ws.Range("I74").Value = 1813.0294
$ws.Range("J74").Value = 6315.4
$ws.Range("K74").Value = 1813.0294
$ws.Range("L74").Value = 6315.4
$ws.Range("M74").Value = -939.0293999999999
$ws.Range("N74").Value = -8063.4

$ws.Range("H77").Value = 2390.2563
$ws.Range("I77").Value = 1813.0294
$ws.Range("J77").Value = 6315.4
$ws.Range("K77").Value = 9065.146999999999
$ws.Range("L77").Value = 31577
$ws.Range("M77").Value = -4697.146999999999
$ws.Range("N77").Value = -40313

$ws.Range("H122").Value = 3053.261
$ws.Range("I122").Value = 2644.2354
$ws.Range("K122").Value = 7932.706200000001
$ws.Range("M122").Value = -5482.706200000001

$ws.Range("H136").Value = 4588.5
$ws.Range("I136").Value = 2537.2727
$ws.Range("K136").Value = 7611.8181
$ws.Range("M136").Value = -5061.8181

$ws.Range("H138").Value = 97000
$ws.Range("J138").Value = 97000
$ws.Range("L138").Value = 97000
$ws.Range("N138").Value = -107280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H39").Value = 9999
$ws.Range("I39").Value = 9999
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 9999
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -9610
$ws.Range("N39").ClearContents()

$ws.Range("H76").Value = 54500
$ws.Range("J76").Value = 54500
$ws.Range("L76").Value = 54500
$ws.Range("N76").Value = -55130

$ws.Range("H79").Value = 54500
$ws.Range("J79").Value = 54500
$ws.Range("L79").Value = 54500
$ws.Range("N79").Value = -56684

$ws.Range("H80").Value = 22540.2
$ws.Range("J80").Value = 15668.875
$ws.Range("L80").Value = 15668.875
$ws.Range("N80").Value = -17664.875

$ws.Range("H83").Value = 22540.2
$ws.Range("J83").Value = 15668.875
$ws.Range("L83").Value = 78344.375
$ws.Range("N83").Value = -88328.375

$ws.Range("H94").Value = 1142.9524
$ws.Range("I94").Value = 939.93335
$ws.Range("J94").Value = 1650.5
$ws.Range("K94").Value = 939.93335
$ws.Range("L94").Value = 1650.5
$ws.Range("M94").Value = -488.93335
$ws.Range("N94").Value = -2552.5

$ws.Range("H134").Value = 5372.533
$ws.Range("I134").Value = 3549.1667
$ws.Range("J134").Value = 12666
$ws.Range("K134").Value = 10647.5001
$ws.Range("L134").Value = 37998
$ws.Range("M134").Value = -8112.500100000001
$ws.Range("N134").Value = -43068

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9150.781000000001
$ws.Range("I31").Value = 3530.4375
$ws.Range("J31").Value = 14771.125
$ws.Range("K31").Value = 3530.4375
$ws.Range("L31").Value = 14771.125
$ws.Range("M31").Value = -3235.4375
$ws.Range("N31").Value = -15361.125

$ws.Range("H34").Value = 9150.781000000001
$ws.Range("I34").Value = 3530.4375
$ws.Range("J34").Value = 14771.125
$ws.Range("K34").Value = 3530.4375
$ws.Range("L34").Value = 14771.125
$ws.Range("M34").Value = -3328.4375
$ws.Range("N34").Value = -15175.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 880.0625
$ws.Range("J5").Value = 1185.4286
$ws.Range("L5").Value = 3556.2858
$ws.Range("N5").Value = -3780.2858

$ws.Range("H70").Value = 10526.5
$ws.Range("I70").Value = 1053
$ws.Range("K70").Value = 3159
$ws.Range("M70").Value = -2844

$ws.Range("H73").Value = 10526.5
$ws.Range("I73").Value = 1053
$ws.Range("K73").Value = 3159
$ws.Range("M73").Value = -2067

$ws.Range("H129").Value = 11113869
$ws.Range("I129").Value = 3613.3333
$ws.Range("J129").Value = 13891433
$ws.Range("K129").Value = 10839.9999
$ws.Range("L129").Value = 41674299
$ws.Range("M129").Value = -5839.999899999999
$ws.Range("N129").Value = -41684299

$ws.Range("H135").Value = 880.0625
$ws.Range("J135").Value = 1185.4286
$ws.Range("L135").Value = 10668.8574
$ws.Range("N135").Value = -15738.8574

$ws.Range("H141").Value = 6879.3335
$ws.Range("I141").Value = 5194.1816
$ws.Range("K141").Value = 15582.5448
$ws.Range("M141").Value = -10402.5448

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4354.778
$ws.Range("I80").Value = 4199.75
$ws.Range("K80").Value = 4199.75
$ws.Range("M80").Value = -3201.75

$ws.Range("H83").Value = 4354.778
$ws.Range("I83").Value = 4199.75
$ws.Range("K83").Value = 20998.75
$ws.Range("M83").Value = -16006.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1081
$ws.Range("I16").Value = 1131.9231
$ws.Range("J16").Value = 419
$ws.Range("K16").Value = 1131.9231
$ws.Range("L16").Value = 419
$ws.Range("M16").Value = -961.9231
$ws.Range("N16").Value = -759

$ws.Range("H22").Value = 4362.5
$ws.Range("I22").Value = 4983.3335
$ws.Range("K22").Value = 4983.3335
$ws.Range("M22").Value = -4688.3335

$ws.Range("H27").Value = 4362.5
$ws.Range("I27").Value = 4983.3335
$ws.Range("K27").Value = 4983.3335
$ws.Range("M27").Value = -4876.3335

$ws.Range("H46").Value = 2926.5
$ws.Range("I46").Value = 1143
$ws.Range("J46").Value = 3818.25
$ws.Range("K46").Value = 1143
$ws.Range("L46").Value = 3818.25
$ws.Range("M46").Value = -955
$ws.Range("N46").Value = -4194.25

$ws.Range("H68").Value = 7100.263
$ws.Range("I68").Value = 3987.25
$ws.Range("J68").Value = 9364.272000000001
$ws.Range("K68").Value = 3987.25
$ws.Range("L68").Value = 9364.272000000001
$ws.Range("M68").Value = -3238.25
$ws.Range("N68").Value = -10862.272

$ws.Range("H71").Value = 7100.263
$ws.Range("I71").Value = 3987.25
$ws.Range("J71").Value = 9364.272000000001
$ws.Range("K71").Value = 19936.25
$ws.Range("L71").Value = 46821.36
$ws.Range("M71").Value = -16192.25
$ws.Range("N71").Value = -54309.36

$ws.Range("H82").Value = 1241.3462
$ws.Range("I82").Value = 799.9286
$ws.Range("J82").Value = 1756.3334
$ws.Range("K82").Value = 799.9286
$ws.Range("L82").Value = 1756.3334
$ws.Range("M82").Value = -438.9286
$ws.Range("N82").Value = -2478.3334

$ws.Range("H85").Value = 1241.3462
$ws.Range("I85").Value = 799.9286
$ws.Range("J85").Value = 1756.3334
$ws.Range("K85").Value = 799.9286
$ws.Range("L85").Value = 1756.3334
$ws.Range("M85").Value = 448.0714
$ws.Range("N85").Value = -4252.3334

$ws.Range("H100").Value = 3924.8125
$ws.Range("J100").Value = 5562.375
$ws.Range("L100").Value = 5562.375
$ws.Range("N100").Value = -6644.375

$ws.Range("H127").Value = 81258.336
$ws.Range("J127").Value = 81258.336
$ws.Range("L127").Value = 81258.336
$ws.Range("N127").Value = -91178.336

$ws.Range("H132").Value = 9249.1875
$ws.Range("I132").Value = 12837.4
$ws.Range("J132").Value = 7618.1816
$ws.Range("K132").Value = 38512.2
$ws.Range("L132").Value = 22854.5448
$ws.Range("M132").Value = -35982.2
$ws.Range("N132").Value = -27914.5448

$ws.Range("H136").Value = 9305.843999999999
$ws.Range("I136").Value = 5778.727
$ws.Range("K136").Value = 17336.181
$ws.Range("M136").Value = -14786.181

$ws.Range("H138").Value = 100000
$ws.Range("J138").Value = 100000
$ws.Range("L138").Value = 100000
$ws.Range("N138").Value = -110280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 59999.5
$ws.Range("J135").Value = 59999.5
$ws.Range("L135").Value = 59999.5
$ws.Range("N135").Value = -70139.5

$ws.Range("H136").Value = 5658.591
$ws.Range("I136").Value = 4394.316
$ws.Range("K136").Value = 13182.948
$ws.Range("M136").Value = -10632.948
